$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-20 Friday" "2024-09-21 Saturday"

Replace-Text "581÷7=" "340÷4="
Replace-Text "263÷9=" "863÷6="
Replace-Text "958÷7=" "226÷2="
Replace-Text "848÷7=" "934÷4="
Replace-Text "440÷2=" "105÷2="

Replace-Text "248÷7=" "256÷6="
Replace-Text "783÷3=" "375÷9="
Replace-Text "637÷9=" "220÷8="
Replace-Text "288÷9=" "450÷2="
Replace-Text "291÷5=" "822÷6="

Replace-Text "821÷5=" "544÷6="
Replace-Text "484÷4=" "925÷2="
Replace-Text "416÷3=" "979÷8="
Replace-Text "433÷7=" "756÷2="
Replace-Text "112÷3=" "853÷9="

Replace-Text "525÷9=" "531÷5="
Replace-Text "332÷7=" "524÷3="
Replace-Text "843÷2=" "395÷7="
Replace-Text "661÷7=" "829÷2="
Replace-Text "722÷7=" "352÷6="

Replace-Text "318÷4=" "463÷3="
Replace-Text "525÷4=" "703÷8="
Replace-Text "201÷6=" "411÷2="
Replace-Text "608÷6=" "222÷2="
Replace-Text "660÷2=" "167÷9="
